$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (15 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1345.4546
$ws.Range("I99").Value = 1345.4546
$ws.Range("K99").Value = 4036.3638
$ws.Range("M99").Value = -2538.3638
$ws.Range("H132").Value = 14755.042
$ws.Range("I132").Value = 2014.6984
$ws.Range("J132").Value = 115085.25
$ws.Range("K132").Value = 6044.0952
$ws.Range("L132").Value = 345255.75
$ws.Range("M132").Value = -3514.0952
$ws.Range("N132").Value = -350315.75
$ws.Range("H138").Value = 2609.2615
$ws.Range("I138").Value = 1367.675
$ws.Range("K138").Value = 4103.025
$ws.Range("M138").Value = 1036.975

# --- Sheet: ARM (21 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2226.103
$ws.Range("I74").Value = 2118.0168
$ws.Range("J74").Value = 2934.6667
$ws.Range("K74").Value = 2118.0168
$ws.Range("L74").Value = 2934.6667
$ws.Range("M74").Value = -1244.0168
$ws.Range("N74").Value = -4682.6667
$ws.Range("H77").Value = 2226.103
$ws.Range("I77").Value = 2118.0168
$ws.Range("J77").Value = 2934.6667
$ws.Range("K77").Value = 10590.084
$ws.Range("L77").Value = 14673.3335
$ws.Range("M77").Value = -6222.083999999999
$ws.Range("N77").Value = -23409.3335
$ws.Range("H132").Value = 8198621.5
$ws.Range("I132").Value = 14287002
$ws.Range("J132").Value = 2724.7307
$ws.Range("K132").Value = 42861006
$ws.Range("L132").Value = 8174.1921
$ws.Range("M132").Value = -42858476
$ws.Range("N132").Value = -13234.1921

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 52992
$ws.Range("J124").Value = 52992
$ws.Range("L124").Value = 52992
$ws.Range("N124").Value = -62812
$ws.Range("H130").Value = 49383.5
$ws.Range("J130").Value = 49383.5
$ws.Range("L130").Value = 49383.5
$ws.Range("N130").Value = -59423.5

# --- Sheet: CRP (36 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2191.97
$ws.Range("I31").Value = 938.4231
$ws.Range("J31").Value = 3549.9792
$ws.Range("K31").Value = 938.4231
$ws.Range("L31").Value = 3549.9792
$ws.Range("M31").Value = -643.4231
$ws.Range("N31").Value = -4139.9792
$ws.Range("H34").Value = 2191.97
$ws.Range("I34").Value = 938.4231
$ws.Range("J34").Value = 3549.9792
$ws.Range("K34").Value = 938.4231
$ws.Range("L34").Value = 3549.9792
$ws.Range("M34").Value = -736.4231
$ws.Range("N34").Value = -3953.9792
$ws.Range("H58").Value = 1765.3793
$ws.Range("I58").Value = 1217.4348
$ws.Range("K58").Value = 1217.4348
$ws.Range("M58").Value = -1014.4348
$ws.Range("H122").Value = 39816.453
$ws.Range("I122").Value = 71674.586
$ws.Range("J122").Value = 1131.5714
$ws.Range("K122").Value = 215023.758
$ws.Range("L122").Value = 3394.7142
$ws.Range("M122").Value = -212573.758
$ws.Range("N122").Value = -8294.7142
$ws.Range("H134").Value = 2229.1428
$ws.Range("I134").Value = 1437.3529
$ws.Range("J134").Value = 3452.818
$ws.Range("K134").Value = 4312.0587
$ws.Range("L134").Value = 10358.454
$ws.Range("M134").Value = -1777.0587
$ws.Range("N134").Value = -15428.454
$ws.Range("H136").Value = 1765.3793
$ws.Range("I136").Value = 1217.4348
$ws.Range("K136").Value = 3652.3044
$ws.Range("M136").Value = -1102.3044

# --- Sheet: CUL (21 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4073.1714
$ws.Range("I5").Value = 9983.272000000001
$ws.Range("J5").Value = 1364.375
$ws.Range("K5").Value = 29949.816
$ws.Range("L5").Value = 4093.125
$ws.Range("M5").Value = -29837.816
$ws.Range("N5").Value = -4317.125
$ws.Range("H122").Value = 3799.5454
$ws.Range("I122").Value = 463.625
$ws.Range("J122").Value = 6939.2354
$ws.Range("K122").Value = 4172.625
$ws.Range("L122").Value = 62453.11859999999
$ws.Range("M122").Value = -1722.625
$ws.Range("N122").Value = -67353.11859999999
$ws.Range("H135").Value = 4073.1714
$ws.Range("I135").Value = 9983.272000000001
$ws.Range("J135").Value = 1364.375
$ws.Range("K135").Value = 89849.448
$ws.Range("L135").Value = 12279.375
$ws.Range("M135").Value = -87314.448
$ws.Range("N135").Value = -17349.375

# --- Sheet: GSM (7 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2384.9285
$ws.Range("I132").Value = 1397.591
$ws.Range("J132").Value = 3471
$ws.Range("K132").Value = 4192.772999999999
$ws.Range("L132").Value = 10413
$ws.Range("M132").Value = -1662.772999999999
$ws.Range("N132").Value = -15473

# --- Sheet: LTW (55 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3040.611
$ws.Range("I7").Value = 2540.4614
$ws.Range("J7").Value = 4341
$ws.Range("K7").Value = 2540.4614
$ws.Range("L7").Value = 4341
$ws.Range("M7").Value = -2428.4614
$ws.Range("N7").Value = -4565
$ws.Range("H61").Value = 2598.923
$ws.Range("I61").Value = 2595.111
$ws.Range("J61").Value = 2607.5
$ws.Range("K61").Value = 2595.111
$ws.Range("L61").Value = 2607.5
$ws.Range("M61").Value = -2393.111
$ws.Range("N61").Value = -3011.5
$ws.Range("H76").Value = 13686.2
$ws.Range("J76").Value = 16390
$ws.Range("L76").Value = 16390
$ws.Range("N76").Value = -17066
$ws.Range("H79").Value = 13686.2
$ws.Range("J79").Value = 16390
$ws.Range("L79").Value = 16390
$ws.Range("N79").Value = -18730
$ws.Range("H81").Value = 32181
$ws.Range("J81").Value = 32181
$ws.Range("L81").Value = 32181
$ws.Range("N81").Value = -34177
$ws.Range("H84").Value = 32181
$ws.Range("J84").Value = 32181
$ws.Range("L84").Value = 96543
$ws.Range("N84").Value = -106527
$ws.Range("H113").Value = 2598.923
$ws.Range("I113").Value = 2595.111
$ws.Range("J113").Value = 2607.5
$ws.Range("K113").Value = 2595.111
$ws.Range("L113").Value = 2607.5
$ws.Range("M113").Value = -425.1109999999999
$ws.Range("N113").Value = -6947.5
$ws.Range("H126").Value = 3040.611
$ws.Range("I126").Value = 2540.4614
$ws.Range("J126").Value = 4341
$ws.Range("K126").Value = 7621.3842
$ws.Range("L126").Value = 13023
$ws.Range("M126").Value = -5151.3842
$ws.Range("N126").Value = -17963
$ws.Range("H127").Value = 50495.332
$ws.Range("J127").Value = 50495.332
$ws.Range("L127").Value = 50495.332
$ws.Range("N127").Value = -60415.332
$ws.Range("H132").Value = 2232.0159
$ws.Range("I132").Value = 1378.3903
$ws.Range("J132").Value = 3822.8635
$ws.Range("K132").Value = 4135.1709
$ws.Range("L132").Value = 11468.5905
$ws.Range("M132").Value = -1605.1709
$ws.Range("N132").Value = -16528.5905

# --- Sheet: WVR (17 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 52828
$ws.Range("J8").Value = 52828
$ws.Range("L8").Value = 52828
$ws.Range("N8").Value = -53108
$ws.Range("H132").Value = 1525.4889
$ws.Range("I132").Value = 1480.614
$ws.Range("J132").Value = 1603
$ws.Range("K132").Value = 4441.842000000001
$ws.Range("L132").Value = 4809
$ws.Range("M132").Value = -1911.842000000001
$ws.Range("N132").Value = -9869
$ws.Range("H136").Value = 733.80554
$ws.Range("I136").Value = 523.2766
$ws.Range("J136").Value = 1129.6
$ws.Range("K136").Value = 1569.8298
$ws.Range("M136").Value = 980.1702
$ws.Range("N136").Value = -8488.799999999999

Write-Output "Applied 180 cell updates across 8 sheets"